# Apply updated odds values to Sheet1 as described by the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new numeric value
$updates = @{
    "M3"  = 1.08
    "N3"  = 8

    "M5"  = 1.08
    "N5"  = 8

    "G6"  = 2.22
    "H6"  = 3.2
    "J6"  = 2.77
    "L6"  = 3.55
    "P6"  = 2.95
    "W6"  = 7.9
    "X6"  = 11.25
    "Y6"  = 8.75
    "AA6" = 18
    "AB6" = 27
    "AD6" = 6.2
    "AM6" = 500
    "AN6" = 4.15
    "AO6" = 11.25
    "AT6" = 2.6
    "AU6" = 6.7
    "AV6" = 55
    "AX6" = 16.5
    "AY6" = 23
    "AZ6" = 80
    "BA6" = 110

    "Q12" = 2.05
    "R12" = 1.75

    "G16"  = 2.3
    "H16"  = 3.25
    "I16"  = 3.1
    "J16"  = 3
    "L16"  = 3.75
    "O16"  = 1.3
    "P16"  = 3.4
    "Z16"  = 21
    "AA16" = 19
    "AG16" = 9.5
    "AI16" = 12
    "AJ16" = 34

    "M17" = 1.03
    "N17" = 15

    "N20"  = 7.5
    "AE20" = 13
    "AG20" = 9.5
    "AQ20" = 41
    "AV20" = 51
    "BB20" = 201

    "H25"  = 3.4
    "I25"  = 3.05
    "R25"  = 2
    "X25"  = 11
    "AD25" = 6.8
    "AE25" = 12.5
    "AG25" = 11.25
    "AH25" = 17.5
    "AJ25" = 40
    "AU25" = 6.8
    "AW25" = 5.2
    "AX25" = 16.5

    "N26" = 5.75
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
